$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '63.636.99'
Set-TextValue 'E2' '  -2.65%  '

Set-TextValue 'D3' '3.334.90'
Set-TextValue 'E3' '  -3.05%  '

Set-TextValue 'E4' '  -0.07%  '

Set-TextValue 'D5' '546.91'
Set-TextValue 'E5' '  -0.57%  '

Set-TextValue 'D6' '171.89'
Set-TextValue 'E6' '  -4.35%  '

Set-TextValue 'E7' '  -4.23%  '

Set-TextValue 'D8' '3.325.48'
Set-TextValue 'E8' '  -3.11%  '

Set-TextValue 'E9' '  -0.08%  '

Set-TextValue 'E10' '  -2.03%  '

Set-TextValue 'E11' '  +0.80%  '

Set-TextValue 'D12' '53.50'
Set-TextValue 'E12' '  +0.19%  '

Set-TextValue 'D13' '0.0000264'
Set-TextValue 'E13' '  -2.11%  '

Set-TextValue 'D14' '8.91'
Set-TextValue 'E14' '  -2.72%  '

Set-TextValue 'D15' '3.877.92'
Set-TextValue 'E15' '  -3.42%  '

Set-TextValue 'B16' 'Chainlink'
Set-TextValue 'C16' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D16' '17.88'
Set-TextValue 'E16' '  -2.26%  '

Set-TextValue 'B17' 'WrappedEther'
Set-TextValue 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '3.332.45'
Set-TextValue 'E17' '  -3.61%  '

Set-TextValue 'E18' '  -3.30%  '

Set-TextValue 'D19' '11.71'
Set-TextValue 'E19' '  -0.87%  '

Set-TextValue 'D20' '63.619.52'
Set-TextValue 'E20' '  -2.87%  '

Set-TextValue 'D21' '0.972'
Set-TextValue 'E21' '  -0.85%  '

Set-TextValue 'D22' '413.65'
Set-TextValue 'E22' '  -0.92%  '

Set-TextValue 'D23' '4.04'
Set-TextValue 'E23' '  +0.53%  '

Set-TextValue 'D24' '4.30'
Set-TextValue 'E24' '  +4.92%  '

Set-TextValue 'D25' '13.65'
Set-TextValue 'E25' '  +12.47%  '

Set-TextValue 'D26' '83.00'
Set-TextValue 'E26' '  -2.08%  '

Set-TextValue 'D27' '10.56'
Set-TextValue 'E27' '  -1.93%  '

Set-TextValue 'E28' '  -4.98%  '

Set-TextValue 'D29' '8.57'
Set-TextValue 'E29' '  -3.74%  '

Set-TextValue 'D30' '29.05'
Set-TextValue 'E30' '  -2.41%  '

Set-TextValue 'D31' '6.37'
Set-TextValue 'E31' '  -2.66%  '

Set-TextValue 'D32' '11.32'
Set-TextValue 'E32' '  -2.66%  '

Set-TextValue 'D33' '574.40'
Set-TextValue 'E33' '  -6.46%  '

Set-TextValue 'E34' '  -2.77%  '

Set-TextValue 'D35' '57.56'
Set-TextValue 'E35' '  -3.28%  '

Set-TextValue 'B36' 'Kaspa'
Set-TextValue 'C36' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D36' '0.147'
Set-TextValue 'E36' '  -0.21%  '

Set-TextValue 'B37' 'Dai'
Set-TextValue 'C37' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D37' '1.00'
Set-TextValue 'E37' '  +0.16%  '

Set-TextValue 'D38' '35.05'
Set-TextValue 'E38' '  -6.01%  '

Set-TextValue 'B39' 'Stacks'
Set-TextValue 'C39' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '3.39'
Set-TextValue 'E39' '  +1.16%  '

Set-TextValue 'B40' 'PEPE'
Set-TextValue 'C40' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D40' '0.0₃0736'
Set-TextValue 'E40' '  -5.39%  '

Set-TextValue 'D41' '0.366'
Set-TextValue 'E41' '  -3.07%  '

Set-TextValue 'D42' '3.141.22'
Set-TextValue 'E42' '  +0.47%  '

Set-TextValue 'E43' '  -0.08%  '

Set-TextValue 'E44' '  +0.16%  '

Set-TextValue 'E45' '  -0.55%  '

Set-TextValue 'D46' '0.0400'
Set-TextValue 'E46' '  -2.20%  '

Set-TextValue 'D47' '2.40'
Set-TextValue 'E47' '  -5.50%  '

Set-TextValue 'D48' '2.59'
Set-TextValue 'E48' '  -4.67%  '

Set-TextValue 'E49' '  -2.84%  '

Set-TextValue 'D50' '132.65'
Set-TextValue 'E50' '  -3.89%  '

Set-TextValue 'D51' '8.02'
Set-TextValue 'E51' '  -3.64%  '
